$wb = $excel.ActiveWorkbook
$papers = $wb.Worksheets.Item("papers")

$papers.Range("A2").Value = "allstadt15"
$papers.Range("B2").Value = ""
$papers.Range("C2").Value = ""
$papers.Range("D2").Value = ""
$papers.Range("E2").Value = ""
$papers.Range("F2").Value = ""
$papers.Range("G2").Value = ""
$papers.Range("H2").Value = ""
$papers.Range("I2").Value = ""
$papers.Range("J2").Value = ""
$papers.Range("K2").Value = ""
$papers.Range("L2").Value = ""
$papers.Range("M2").Value = ""
$papers.Range("N2").Value = ""
$papers.Range("O2").Value = ""
$papers.Range("P2").Value = ""
$papers.Range("Q2").Value = ""
$papers.Range("R2").Value = ""
$papers.Range("S2").Value = ""
$papers.Range("T2").Value = ""
$papers.Range("U2").Value = ""
$papers.Range("V2").Value = ""
$papers.Range("W2").Value = ""
$papers.Range("A3").Value = "augspurger13"
$papers.Range("A4").Value = "ault13"
$papers.Range("A5").Value = "ault15"
$papers.Range("A6").Value = "basler12"
$papers.Range("A7").Value = "caradonna16"
$papers.Range("A8").Value = "fu12"
$papers.Range("A9").Value = "gu08"
$papers.Range("A10").Value = "hufkens12"
$papers.Range("A11").Value = "kodra11"
$papers.Range("A12").Value = "koehler12"
$papers.Range("A13").Value = "kollas14"
$papers.Range("A14").Value = "korner10"
$papers.Range("A15").Value = "korner16"
$papers.Range("A16").Value = "lenz13"
$papers.Range("B16").Value = "Lenz, A."
$papers.Range("C16").Value = 2013
$papers.Range("D16").Value = "New Phytologist"
$papers.Range("E16").Value = "European deciduous trees exhibit similar safety margins against damage by spring freeze events along elevational gradients"
$papers.Range("F16").Value = "1-10"
$papers.Range("G16").Value = "10.1111/nph.12452"
$papers.Range("H16").Value = "observational"
$papers.Range("I16").Value = "Sorbus aucuparia, S. aria, Acer psudoplatanus, Laburnum alpinum, Prunus avium, Fagus sylvatica, Fraxinus excelsior"
$papers.Range("J16").Value = -35
$papers.Range("K16").Value = 4
$papers.Range("L16").Value = "chilling and forcing - 3K/hr"
$papers.Range("M16").Value = "no"
$papers.Range("N16").Value = "both"
$papers.Range("O16").Value = "yes"
$papers.Range("P16").Value = "no"
$papers.Range("Q16").Value = "chilling and forcing"
$papers.Range("R16").Value = "yes"
$papers.Range("S16").Value = "yes - freezing resistence"
$papers.Range("T16").Value = "no"
$papers.Range("U16").Value = "no"
$papers.Range("V16").Value = "no"
$papers.Range("W16").Value = "elevational gradient and safety margins against freezing damage"
$papers.Range("A21").Value = "polgar11"
$papers.Range("A22").Value = "polgar14"
$papers.Range("A26").Value = "schwartz10"
$papers.Range("A27").Value = "schwartz13"
$papers.Range("A28").Value = "schwartz90"
$papers.Range("A29").Value = "vavrus06"
$papers.Range("A30").Value = "vitasse13"
$papers.Range("A31").Value = "vitasse13"
$papers.Range("A32").Value = "vitasse14"

$papers.Range("A1:W1").AutoFilter()

$searches = $wb.Worksheets.Item("searches")
$searches.Range("D5").Select()

$papers.Activate()
$papers.Range("E8").Select()
